# Refresh derived market-price / leve-profit figures (scheduled runner update)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 801.6
$ws.Range("I19").Value = 933.3333
$ws.Range("J19").Value = 693.8182
$ws.Range("K19").Value = 933.3333
$ws.Range("L19").Value = 693.8182
$ws.Range("M19").Value = -758.3333
$ws.Range("N19").Value = -1043.8182
$ws.Range("H32").Value = 883.1818
$ws.Range("I32").Value = 655.8570999999999
$ws.Range("J32").Value = 989.26666
$ws.Range("K32").Value = 655.8570999999999
$ws.Range("L32").Value = 989.26666
$ws.Range("M32").Value = -329.8570999999999
$ws.Range("N32").Value = -1641.26666
$ws.Range("H33").Value = 382.0811
$ws.Range("I33").Value = 115
$ws.Range("K33").Value = 115
$ws.Range("M33").Value = 114
$ws.Range("H112").Value = 1111.08
$ws.Range("J112").Value = 1138.6957
$ws.Range("L112").Value = 3416.0871
$ws.Range("N112").Value = -5632.0871
$ws.Range("H134").Value = 46880
$ws.Range("J134").Value = 46880
$ws.Range("L134").Value = 46880
$ws.Range("N134").Value = -57020
$ws.Range("H137").Value = 43406.918
$ws.Range("I137").Value = 92356
$ws.Range("J137").Value = 1988.4615
$ws.Range("K137").Value = 277068
$ws.Range("L137").Value = 5965.3845
$ws.Range("M137").Value = -274518
$ws.Range("N137").Value = -11065.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1732.6342
$ws.Range("I61").Value = 1624
$ws.Range("J61").Value = 2514.8
$ws.Range("K61").Value = 1624
$ws.Range("L61").Value = 2514.8
$ws.Range("M61").Value = -1412
$ws.Range("N61").Value = -2938.8
$ws.Range("H136").Value = 1732.6342
$ws.Range("I136").Value = 1624
$ws.Range("J136").Value = 2514.8
$ws.Range("K136").Value = 4872
$ws.Range("L136").Value = 7544.400000000001
$ws.Range("M136").Value = -2322
$ws.Range("N136").Value = -12644.4
$ws.Range("H139").Value = 61857.332
$ws.Range("J139").Value = 61857.332
$ws.Range("L139").Value = 61857.332
$ws.Range("N139").Value = -72137.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 992.9091
$ws.Range("I107").Value = 902.44446
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 902.44446
$ws.Range("L107").Value = 1400
$ws.Range("M107").Value = 1017.55554
$ws.Range("N107").Value = -5240
$ws.Range("H134").Value = 33402418
$ws.Range("I134").Value = 45455300
$ws.Range("J134").Value = 256993.5
$ws.Range("K134").Value = 136365900
$ws.Range("L134").Value = 770980.5
$ws.Range("M134").Value = -136363365
$ws.Range("N134").Value = -776050.5
$ws.Range("H138").Value = 53115.293
$ws.Range("I138").Value = 10000
$ws.Range("J138").Value = 55810
$ws.Range("K138").Value = 10000
$ws.Range("L138").Value = 55810
$ws.Range("M138").Value = -4860
$ws.Range("N138").Value = -66090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24135.152
$ws.Range("I31").Value = 32338.031
$ws.Range("J31").Value = 5385.7144
$ws.Range("K31").Value = 32338.031
$ws.Range("L31").Value = 5385.7144
$ws.Range("M31").Value = -32043.031
$ws.Range("N31").Value = -5975.7144
$ws.Range("H34").Value = 24135.152
$ws.Range("I34").Value = 32338.031
$ws.Range("J34").Value = 5385.7144
$ws.Range("K34").Value = 32338.031
$ws.Range("L34").Value = 5385.7144
$ws.Range("M34").Value = -32136.031
$ws.Range("N34").Value = -5789.7144
$ws.Range("H100").Value = 66245
$ws.Range("J100").Value = 66245
$ws.Range("L100").Value = 66245
$ws.Range("N100").Value = -68409

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2920.6667
$ws.Range("I3").Value = 2190
$ws.Range("J3").Value = 3943.6
$ws.Range("K3").Value = 6570
$ws.Range("L3").Value = 11830.8
$ws.Range("M3").Value = -6458
$ws.Range("N3").Value = -12054.8
$ws.Range("H56").Value = 3419.682
$ws.Range("I56").Value = 3419.682
$ws.Range("K56").Value = 3419.682
$ws.Range("M56").Value = -2889.682
$ws.Range("H68").Value = 11708.556
$ws.Range("I68").Value = 20314.2
$ws.Range("J68").Value = 951.5
$ws.Range("K68").Value = 60942.60000000001
$ws.Range("L68").Value = 2854.5
$ws.Range("M68").Value = -60131.60000000001
$ws.Range("N68").Value = -4476.5
$ws.Range("H71").Value = 11708.556
$ws.Range("I71").Value = 20314.2
$ws.Range("J71").Value = 951.5
$ws.Range("K71").Value = 182827.8
$ws.Range("L71").Value = 8563.5
$ws.Range("M71").Value = -178771.8
$ws.Range("N71").Value = -16675.5
$ws.Range("H107").Value = 326.66666
$ws.Range("I107").Value = 340
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 1020
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 900
$ws.Range("N107").Value = -4740
$ws.Range("H113").Value = 763.8
$ws.Range("J113").Value = 763.8
$ws.Range("L113").Value = 2291.4
$ws.Range("N113").Value = -6631.4
$ws.Range("H122").Value = 756.0833
$ws.Range("J122").Value = 896.55554
$ws.Range("L122").Value = 8068.99986
$ws.Range("N122").Value = -12968.99986
$ws.Range("H131").Value = 948.8
$ws.Range("J131").Value = 1006.88464
$ws.Range("L131").Value = 3020.65392
$ws.Range("N131").Value = -13100.65392
$ws.Range("H132").Value = 513
$ws.Range("I132").Value = 513
$ws.Range("K132").Value = 4617
$ws.Range("M132").Value = -2087

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 80033
$ws.Range("J52").Value = 80033
$ws.Range("L52").Value = 80033
$ws.Range("N52").Value = -80551
$ws.Range("H132").Value = 29353.816
$ws.Range("I132").Value = 1825.7894
$ws.Range("J132").Value = 56881.844
$ws.Range("K132").Value = 5477.3682
$ws.Range("L132").Value = 170645.532
$ws.Range("M132").Value = -2947.3682
$ws.Range("N132").Value = -175705.532

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1043.4445
$ws.Range("I46").Value = 800.25
$ws.Range("J46").Value = 1238
$ws.Range("K46").Value = 800.25
$ws.Range("L46").Value = 1238
$ws.Range("M46").Value = -612.25
$ws.Range("N46").Value = -1614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360
